$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 2019
$ws.Range("C2").Value = 2213
$ws.Range("D2").Value = 106
$ws.Range("E2").Value = 2107
$ws.Range("F2").Value = 10.71789686552073
$ws.Range("G2").Value = 95.21012200632626
$ws.Range("H2").Value = 4.789877993673746

# Row 3 - 2020
$ws.Range("C3").Value = 2820
$ws.Range("D3").Value = 137
$ws.Range("E3").Value = 2683
$ws.Range("F3").Value = 6.190691369182106
$ws.Range("G3").Value = 95.1418439716312
$ws.Range("H3").Value = 4.858156028368795

# Row 4 - 2021
$ws.Range("C4").Value = 2461
$ws.Range("D4").Value = 178
$ws.Range("E4").Value = 2283
$ws.Range("F4").Value = 6.312056737588652
$ws.Range("G4").Value = 92.76716781796019
$ws.Range("H4").Value = 7.232832182039821

# Row 5 - 2022
$ws.Range("C5").Value = 2749
$ws.Range("D5").Value = 314
$ws.Range("E5").Value = 2435
$ws.Range("F5").Value = 12.75904104022755
$ws.Range("G5").Value = 88.57766460531103
$ws.Range("H5").Value = 11.42233539468898

# Row 6 - 2023
$ws.Range("C6").Value = 2734
$ws.Range("D6").Value = 317
$ws.Range("E6").Value = 2417
$ws.Range("F6").Value = 11.53146598763187
$ws.Range("G6").Value = 88.40526700804682
$ws.Range("H6").Value = 11.59473299195318

# Row 7 - 2024
$ws.Range("C7").Value = 2438
$ws.Range("D7").Value = 301
$ws.Range("E7").Value = 2137
$ws.Range("F7").Value = 11.00950987564009
$ws.Range("G7").Value = 87.6538146021329
$ws.Range("H7").Value = 12.3461853978671

# Row 8 - 2025
$ws.Range("C8").Value = 761
$ws.Range("D8").Value = 137
$ws.Range("E8").Value = 624
$ws.Range("F8").Value = 5.619360131255127
$ws.Range("G8").Value = 81.99737187910644
$ws.Range("H8").Value = 18.00262812089356
